$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NumSecu (colonne B) pour les patients existants
$ws.Range("B2").Value = 181025900000001
$ws.Range("B3").Value = 280055900000002
$ws.Range("B4").Value = 111115900000003
$ws.Range("B5").Value = 164126200000004
$ws.Range("B6").Value = 297110300000005
$ws.Range("B7").Value = 176076200000006
$ws.Range("B8").Value = 286015900000007
$ws.Range("B9").Value = 202098100000008
$ws.Range("B10").Value = 194037600000009
$ws.Range("B11").Value = 153025900000010
$ws.Range("B12").Value = 196016200000011
$ws.Range("B13").Value = 182045900000012
$ws.Range("B14").Value = 207105900000013
$ws.Range("B15").Value = 280094400000014
$ws.Range("B16").Value = 178035900000015
$ws.Range("B17").Value = 281086200000016
$ws.Range("B18").Value = 208085900000017
$ws.Range("B19").Value = 113035900000018
$ws.Range("B20").Value = 192105900000019
$ws.Range("B21").Value = 166025900000020
$ws.Range("B22").Value = 271123700000021
$ws.Range("B23").Value = 200020100000022
$ws.Range("B24").Value = 179035900000023
$ws.Range("B25").Value = 146065400000024
$ws.Range("B26").Value = 102045900000025
$ws.Range("B27").Value = 254065900000056
$ws.Range("B28").Value = 296096200000027
$ws.Range("B29").Value = 201015900000028
$ws.Range("B30").Value = 188065900000029
$ws.Range("B31").Value = 174075900000030
$ws.Range("B32").Value = 199126200000031
$ws.Range("B33").Value = 175115900000032
$ws.Range("B34").Value = 163055900000033
$ws.Range("B35").Value = 160095900000034
$ws.Range("B36").Value = 262015900000035
$ws.Range("B37").Value = 183047600000036
$ws.Range("B38").Value = 119065900000037

# Textes (dans l ordre exact de saisie, pour reconstituer la table des chaines partagees)
$ws.Range("J28").Value = '09 87 52 63 33'
$ws.Range("J29").Value = '06 45 65 98 98'
$ws.Range("J30").Value = '03 20 50 40 40'
$ws.Range("J31").Value = '03 21 30 30 21 '
$ws.Range("J32").Value = '03 20 30 20 03'
$ws.Range("J33").Value = '03 21 21 30 21'
$ws.Range("J34").Value = '03 32 35 84 84'
$ws.Range("J35").Value = '03 23 32 10 51 '
$ws.Range("J36").Value = '03 21 21 87 89'
$ws.Range("J37").Value = '09 02 22 08 34'
$ws.Range("J39").Value = '03 20 20 20 20'
$ws.Range("J41").Value = '03 21 11 22 33'
$ws.Range("J42").Value = '03 21 88 88 99'
$ws.Range("J43").Value = '03 20 22 55 88'
$ws.Range("J47").Value = '09 10 56 98 33'
$ws.Range("J48").Value = '09 80 51 23 44'
$ws.Range("J49").Value = '09 87 65 32 15'
$ws.Range("J50").Value = '09 78 98 78 98'
$ws.Range("J51").Value = '09 23 65 98 15'
$ws.Range("J52").Value = '09 88 98 98 65 '
$ws.Range("J53").Value = '03 20 21 21 60'
$ws.Range("J55").Value = '09 89 63 36 65'
$ws.Range("J54").Value = '09 56 65 98 12'
$ws.Range("J56").Value = '03 23 23 84 89'
$ws.Range("J57").Value = '03 20 20 99 99'
$ws.Range("J58").Value = '03 21 87 87 87'
$ws.Range("J59").Value = '03 56 56 45 45'
$ws.Range("J60").Value = '03 21 65 98 89 '
$ws.Range("J61").Value = '03 56 89 54 54'
$ws.Range("J62").Value = '09 87 89 32 23'
$ws.Range("J63").Value = '09 87 87 89 98'
$ws.Range("J64").Value = '09 23 56 65 56'
$ws.Range("J65").Value = '03 23 65 65 23'
$ws.Range("J66").Value = '03 21 21 23 88'
$ws.Range("J67").Value = '03 21 66 44 77'
$ws.Range("J68").Value = '03 55 55 88 88'
$ws.Range("J69").Value = '03 23 98 89 77'
$ws.Range("J70").Value = '03 45 65 56 65'
$ws.Range("J71").Value = '03 21 11 44 77'
$ws.Range("G39").Value = '54 Rue des oiseaux'
$ws.Range("G41").Value = '1b place de l''église'
$ws.Range("G42").Value = '37 Rue de Tourcoing'
$ws.Range("G43").Value = '12 Rue Gambetta'
$ws.Range("G47").Value = '23 Place de la monnaie'
$ws.Range("G48").Value = '32 Faubourg Jean Jaurès'
$ws.Range("G49").Value = '135 Rue Nationale'
$ws.Range("G50").Value = '77  Rue Jules Guesde'
$ws.Range("G51").Value = '26 Rue Solférino'
$ws.Range("G52").Value = '23 Grand Place'
$ws.Range("G53").Value = '65 Rue de la soif'
$ws.Range("G54").Value = '32 Rue de la soif'
$ws.Range("G55").Value = '3 Rue Solférino, Apt25'
$ws.Range("G56").Value = '32 Rue de Gand'
$ws.Range("G57").Value = '44 Rue de Gand'
$ws.Range("G58").Value = '23 Rue Nationale'
$ws.Range("G59").Value = '192 Rue Nationale'
$ws.Range("G60").Value = '70 Rue du Général de Gaulles'
$ws.Range("G61").Value = '25 Rue du Général Leclerc'
$ws.Range("G63").Value = '19 Rue de la monnaie'
$ws.Range("G62").Value = '31 Rue de l''école'
$ws.Range("G64").Value = '65 Rue de Roubaix'
$ws.Range("G65").Value = '137 Rue de Tourcoing'
$ws.Range("G66").Value = '26 Allée des sages'
$ws.Range("G67").Value = '25 Boulevard Excellemans'
$ws.Range("G68").Value = '30 Rue Solférino'
$ws.Range("G69").Value = '91 Rue Nationale, apt24'
$ws.Range("G70").Value = '149 Rue du Général de Gaulles'
$ws.Range("G71").Value = '2 Rue de la Monnaie'
$ws.Range("C39").Value = 'LEROY'
$ws.Range("D39").Value = 'Timothé'
$ws.Range("C40").Value = 'DURAND'
$ws.Range("D40").Value = 'Laura'
$ws.Range("C41").Value = 'MENDEZ'
$ws.Range("D41").Value = 'Alicia'
$ws.Range("C45").Value = 'BARGUI'
$ws.Range("D45").Value = 'Mohamed'
$ws.Range("C43").Value = 'N''GUYEN'
$ws.Range("D43").Value = 'Tao'
$ws.Range("D42").Value = 'Virginie'
$ws.Range("C42").Value = 'DETEZ'
$ws.Range("C44").Value = 'DEWAHIE'
$ws.Range("D46").Value = 'Téo'
$ws.Range("C46").Value = 'PETIT'
$ws.Range("D47").Value = 'Carine'
$ws.Range("C48").Value = 'SALUN'
$ws.Range("D48").Value = 'Stéphanie'
$ws.Range("C49").Value = 'SUAREZ'
$ws.Range("D49").Value = 'Arthuro'
$ws.Range("C50").Value = 'VAN BERGUER'
$ws.Range("D50").Value = 'Fanny'
$ws.Range("C51").Value = 'HERTEM'
$ws.Range("D51").Value = 'Téa'
$ws.Range("C52").Value = 'ALDON'
$ws.Range("D53").Value = 'Frédéric'
$ws.Range("D54").Value = 'Marion'
$ws.Range("D55").Value = 'Thibault'
$ws.Range("D57").Value = 'Marie'
$ws.Range("D58").Value = 'Corentin'
$ws.Range("D60").Value = 'Albert'
$ws.Range("D62").Value = 'Stéphane'
$ws.Range("D63").Value = 'Camille'
$ws.Range("D64").Value = 'Alice'
$ws.Range("D65").Value = 'Paul'
$ws.Range("D66").Value = 'Gabin'
$ws.Range("C53").Value = 'TOURU'
$ws.Range("C54").Value = 'BONDRIT'
$ws.Range("C55").Value = 'POURTOI'
$ws.Range("C56").Value = 'MULLIER'
$ws.Range("C57").Value = 'FOURNIER'
$ws.Range("C58").Value = 'SAGE'
$ws.Range("C59").Value = 'PARRAIN'
$ws.Range("C60").Value = 'LUCAS'
$ws.Range("C61").Value = 'MONDON'
$ws.Range("C62").Value = 'CUEILLE'
$ws.Range("C63").Value = 'COGNARD'
$ws.Range("C64").Value = 'MORLIERE'
$ws.Range("C65").Value = 'TREUILLARD'
$ws.Range("C66").Value = 'BADRAN'
$ws.Range("C67").Value = 'TRUSSANT'
$ws.Range("C68").Value = 'GIBOURI'
$ws.Range("C69").Value = 'TRANCHANT'
$ws.Range("C70").Value = 'FILLOT'
$ws.Range("C71").Value = 'DELAUME'

# Completion des nouvelles lignes (39 a 71)
$ws.Range("F2").Copy()
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = 196055900000038
$ws.Range("C39").Value = 'LEROY'
$ws.Range("D39").Value = 'Timothé'
$ws.Range("E39").Value = 'h'
$ws.Range("F39").Value = 35206
$ws.Range("F39").PasteSpecial(-4122)
$ws.Range("G39").Value = '54 Rue des oiseaux'
$ws.Range("H39").Value = 59350
$ws.Range("I39").Value = 'MARCQ EN BAROEUL'
$ws.Range("J39").Value = '03 20 20 20 20'
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = 298015900000039
$ws.Range("C40").Value = 'DURAND'
$ws.Range("D40").Value = 'Laura'
$ws.Range("E40").Value = 'f'
$ws.Range("F40").Value = 35807
$ws.Range("F40").PasteSpecial(-4122)
$ws.Range("G40").Value = '54 Rue des oiseaux'
$ws.Range("H40").Value = 59350
$ws.Range("I40").Value = 'MARCQ EN BAROEUL'
$ws.Range("J40").Value = '03 20 20 20 20'
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = 251026200000040
$ws.Range("C41").Value = 'MENDEZ'
$ws.Range("D41").Value = 'Alicia'
$ws.Range("E41").Value = 'f'
$ws.Range("F41").Value = 18685
$ws.Range("F41").PasteSpecial(-4122)
$ws.Range("G41").Value = '1b place de l''église'
$ws.Range("H41").Value = 59350
$ws.Range("I41").Value = 'MARCQ EN BAROEUL'
$ws.Range("J41").Value = '03 21 11 22 33'
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = 256055900000041
$ws.Range("C42").Value = 'DETEZ'
$ws.Range("D42").Value = 'Virginie'
$ws.Range("E42").Value = 'f'
$ws.Range("F42").Value = 20594
$ws.Range("F42").PasteSpecial(-4122)
$ws.Range("G42").Value = '37 Rue de Tourcoing'
$ws.Range("H42").Value = 59350
$ws.Range("I42").Value = 'MARCQ EN BAROEUL'
$ws.Range("J42").Value = '03 21 88 88 99'
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = 172017500000042
$ws.Range("C43").Value = 'N''GUYEN'
$ws.Range("D43").Value = 'Tao'
$ws.Range("E43").Value = 'h'
$ws.Range("F43").Value = 26299
$ws.Range("F43").PasteSpecial(-4122)
$ws.Range("G43").Value = '12 Rue Gambetta'
$ws.Range("H43").Value = 59350
$ws.Range("I43").Value = 'MARCQ EN BAROEUL'
$ws.Range("J43").Value = '03 20 22 55 88'
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = 275045900000043
$ws.Range("C44").Value = 'DEWAHIE'
$ws.Range("D44").Value = 'Hélène'
$ws.Range("E44").Value = 'f'
$ws.Range("F44").Value = 27507
$ws.Range("F44").PasteSpecial(-4122)
$ws.Range("G44").Value = '12 Rue Gambetta'
$ws.Range("H44").Value = 59350
$ws.Range("I44").Value = 'MARCQ EN BAROEUL'
$ws.Range("J44").Value = '03 20 22 55 88'
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = 199115900000044
$ws.Range("C45").Value = 'BARGUI'
$ws.Range("D45").Value = 'Mohamed'
$ws.Range("E45").Value = 'h'
$ws.Range("F45").Value = 36474
$ws.Range("F45").PasteSpecial(-4122)
$ws.Range("G45").Value = '12 Rue Gambetta'
$ws.Range("H45").Value = 59350
$ws.Range("I45").Value = 'MARCQ EN BAROEUL'
$ws.Range("J45").Value = '03 20 22 55 88'
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = 103035900000045
$ws.Range("C46").Value = 'PETIT'
$ws.Range("D46").Value = 'Téo'
$ws.Range("E46").Value = 'h'
$ws.Range("F46").Value = 37703
$ws.Range("F46").PasteSpecial(-4122)
$ws.Range("G46").Value = '12 Rue Gambetta'
$ws.Range("H46").Value = 59350
$ws.Range("I46").Value = 'MARCQ EN BAROEUL'
$ws.Range("J46").Value = '03 20 22 55 88'
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = 264126200000046
$ws.Range("C47").Value = 'MARTIN'
$ws.Range("D47").Value = 'Carine'
$ws.Range("E47").Value = 'f'
$ws.Range("F47").Value = 23721
$ws.Range("F47").PasteSpecial(-4122)
$ws.Range("G47").Value = '23 Place de la monnaie'
$ws.Range("H47").Value = 59000
$ws.Range("I47").Value = 'LILLE'
$ws.Range("J47").Value = '09 10 56 98 33'
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = 286015900000047
$ws.Range("C48").Value = 'SALUN'
$ws.Range("D48").Value = 'Stéphanie'
$ws.Range("E48").Value = 'f'
$ws.Range("F48").Value = 31423
$ws.Range("F48").PasteSpecial(-4122)
$ws.Range("G48").Value = '32 Faubourg Jean Jaurès'
$ws.Range("H48").Value = 59000
$ws.Range("I48").Value = 'LILLE'
$ws.Range("J48").Value = '09 80 51 23 44'
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = 192115900000048
$ws.Range("C49").Value = 'SUAREZ'
$ws.Range("D49").Value = 'Arthuro'
$ws.Range("E49").Value = 'h'
$ws.Range("F49").Value = 33931
$ws.Range("F49").PasteSpecial(-4122)
$ws.Range("G49").Value = '135 Rue Nationale'
$ws.Range("H49").Value = 59000
$ws.Range("I49").Value = 'LILLE'
$ws.Range("J49").Value = '09 87 65 32 15'
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = 295065900000049
$ws.Range("C50").Value = 'VAN BERGUER'
$ws.Range("D50").Value = 'Fanny'
$ws.Range("E50").Value = 'f'
$ws.Range("F50").Value = 34880
$ws.Range("F50").PasteSpecial(-4122)
$ws.Range("G50").Value = '77  Rue Jules Guesde'
$ws.Range("H50").Value = 59000
$ws.Range("I50").Value = 'LILLE'
$ws.Range("J50").Value = '09 78 98 78 98'
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = 294105900000050
$ws.Range("C51").Value = 'HERTEM'
$ws.Range("D51").Value = 'Téa'
$ws.Range("E51").Value = 'f'
$ws.Range("F51").Value = 34612
$ws.Range("F51").PasteSpecial(-4122)
$ws.Range("G51").Value = '26 Rue Solférino'
$ws.Range("H51").Value = 59000
$ws.Range("I51").Value = 'LILLE'
$ws.Range("J51").Value = '09 23 65 98 15'
$ws.Range("A52").Value = 51
$ws.Range("B52").Value = 277085900000051
$ws.Range("C52").Value = 'ALDON'
$ws.Range("D52").Value = 'Elodie'
$ws.Range("E52").Value = 'f'
$ws.Range("F52").Value = 28350
$ws.Range("F52").PasteSpecial(-4122)
$ws.Range("G52").Value = '23 Grand Place'
$ws.Range("H52").Value = 59000
$ws.Range("I52").Value = 'LILLE'
$ws.Range("J52").Value = '09 88 98 98 65 '
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = 182071200000052
$ws.Range("C53").Value = 'TOURU'
$ws.Range("D53").Value = 'Frédéric'
$ws.Range("E53").Value = 'h'
$ws.Range("F53").Value = 30163
$ws.Range("F53").PasteSpecial(-4122)
$ws.Range("G53").Value = '65 Rue de la soif'
$ws.Range("H53").Value = 59000
$ws.Range("I53").Value = 'LILLE'
$ws.Range("J53").Value = '03 20 21 21 60'
$ws.Range("A54").Value = 53
$ws.Range("B54").Value = 258055900000053
$ws.Range("C54").Value = 'BONDRIT'
$ws.Range("D54").Value = 'Marion'
$ws.Range("E54").Value = 'f'
$ws.Range("F54").Value = 21330
$ws.Range("F54").PasteSpecial(-4122)
$ws.Range("G54").Value = '32 Rue de la soif'
$ws.Range("H54").Value = 59000
$ws.Range("I54").Value = 'LILLE'
$ws.Range("J54").Value = '09 56 65 98 12'
$ws.Range("A55").Value = 54
$ws.Range("B55").Value = 183035900000054
$ws.Range("C55").Value = 'POURTOI'
$ws.Range("D55").Value = 'Thibault'
$ws.Range("E55").Value = 'h'
$ws.Range("F55").Value = 30387
$ws.Range("F55").PasteSpecial(-4122)
$ws.Range("G55").Value = '3 Rue Solférino, Apt25'
$ws.Range("H55").Value = 59000
$ws.Range("I55").Value = 'LILLE'
$ws.Range("J55").Value = '09 89 63 36 65'
$ws.Range("A56").Value = 55
$ws.Range("B56").Value = 100026200000055
$ws.Range("C56").Value = 'MULLIER'
$ws.Range("D56").Value = 'Romain'
$ws.Range("E56").Value = 'h'
$ws.Range("F56").Value = 36585
$ws.Range("F56").PasteSpecial(-4122)
$ws.Range("G56").Value = '32 Rue de Gand'
$ws.Range("H56").Value = 59000
$ws.Range("I56").Value = 'LILLE'
$ws.Range("J56").Value = '03 23 23 84 89'
$ws.Range("A57").Value = 56
$ws.Range("B57").Value = 276095900000056
$ws.Range("C57").Value = 'FOURNIER'
$ws.Range("D57").Value = 'Marie'
$ws.Range("E57").Value = 'f'
$ws.Range("F57").Value = 28017
$ws.Range("F57").PasteSpecial(-4122)
$ws.Range("G57").Value = '44 Rue de Gand'
$ws.Range("H57").Value = 59000
$ws.Range("I57").Value = 'LILLE'
$ws.Range("J57").Value = '03 20 20 99 99'
$ws.Range("A58").Value = 57
$ws.Range("B58").Value = 192065900000057
$ws.Range("C58").Value = 'SAGE'
$ws.Range("D58").Value = 'Corentin'
$ws.Range("E58").Value = 'h'
$ws.Range("F58").Value = 33781
$ws.Range("F58").PasteSpecial(-4122)
$ws.Range("G58").Value = '23 Rue Nationale'
$ws.Range("H58").Value = 59000
$ws.Range("I58").Value = 'LILLE'
$ws.Range("J58").Value = '03 21 87 87 87'
$ws.Range("A59").Value = 58
$ws.Range("B59").Value = 262075900000058
$ws.Range("C59").Value = 'PARRAIN'
$ws.Range("D59").Value = 'Claire'
$ws.Range("E59").Value = 'f'
$ws.Range("F59").Value = 22841
$ws.Range("F59").PasteSpecial(-4122)
$ws.Range("G59").Value = '192 Rue Nationale'
$ws.Range("H59").Value = 59000
$ws.Range("I59").Value = 'LILLE'
$ws.Range("J59").Value = '03 56 56 45 45'
$ws.Range("A60").Value = 59
$ws.Range("B60").Value = 159065900000059
$ws.Range("C60").Value = 'LUCAS'
$ws.Range("D60").Value = 'Albert'
$ws.Range("E60").Value = 'h'
$ws.Range("F60").Value = 21703
$ws.Range("F60").PasteSpecial(-4122)
$ws.Range("G60").Value = '70 Rue du Général de Gaulles'
$ws.Range("H60").Value = 59000
$ws.Range("I60").Value = 'LILLE'
$ws.Range("J60").Value = '03 21 65 98 89 '
$ws.Range("A61").Value = 60
$ws.Range("B61").Value = 154035900000060
$ws.Range("C61").Value = 'MONDON'
$ws.Range("D61").Value = 'Arthur'
$ws.Range("E61").Value = 'h'
$ws.Range("F61").Value = 34415
$ws.Range("F61").PasteSpecial(-4122)
$ws.Range("G61").Value = '25 Rue du Général Leclerc'
$ws.Range("H61").Value = 59000
$ws.Range("I61").Value = 'LILLE'
$ws.Range("J61").Value = '03 56 89 54 54'
$ws.Range("A62").Value = 61
$ws.Range("B62").Value = 185015900000061
$ws.Range("C62").Value = 'CUEILLE'
$ws.Range("D62").Value = 'Stéphane'
$ws.Range("E62").Value = 'h'
$ws.Range("F62").Value = 31064
$ws.Range("F62").PasteSpecial(-4122)
$ws.Range("G62").Value = '31 Rue de l''école'
$ws.Range("H62").Value = 59000
$ws.Range("I62").Value = 'LILLE'
$ws.Range("J62").Value = '09 87 89 32 23'
$ws.Range("A63").Value = 62
$ws.Range("B63").Value = 274065900000062
$ws.Range("C63").Value = 'COGNARD'
$ws.Range("D63").Value = 'Camille'
$ws.Range("E63").Value = 'f'
$ws.Range("F63").Value = 27199
$ws.Range("F63").PasteSpecial(-4122)
$ws.Range("G63").Value = '19 Rue de la monnaie'
$ws.Range("H63").Value = 59000
$ws.Range("I63").Value = 'LILLE'
$ws.Range("J63").Value = '09 87 87 89 98'
$ws.Range("A64").Value = 63
$ws.Range("B64").Value = 298105900000063
$ws.Range("C64").Value = 'MORLIERE'
$ws.Range("D64").Value = 'Alice'
$ws.Range("E64").Value = 'f'
$ws.Range("F64").Value = 36078
$ws.Range("F64").PasteSpecial(-4122)
$ws.Range("G64").Value = '65 Rue de Roubaix'
$ws.Range("H64").Value = 59000
$ws.Range("I64").Value = 'LILLE'
$ws.Range("J64").Value = '09 23 56 65 56'
$ws.Range("A65").Value = 64
$ws.Range("B65").Value = 191030500000064
$ws.Range("C65").Value = 'TREUILLARD'
$ws.Range("D65").Value = 'Paul'
$ws.Range("E65").Value = 'h'
$ws.Range("F65").Value = 33324
$ws.Range("F65").PasteSpecial(-4122)
$ws.Range("G65").Value = '137 Rue de Tourcoing'
$ws.Range("H65").Value = 59000
$ws.Range("I65").Value = 'LILLE'
$ws.Range("J65").Value = '03 23 65 65 23'
$ws.Range("A66").Value = 65
$ws.Range("B66").Value = 183065900000065
$ws.Range("C66").Value = 'BADRAN'
$ws.Range("D66").Value = 'Gabin'
$ws.Range("E66").Value = 'h'
$ws.Range("F66").Value = 30497
$ws.Range("F66").PasteSpecial(-4122)
$ws.Range("G66").Value = '26 Allée des sages'
$ws.Range("H66").Value = 59000
$ws.Range("I66").Value = 'LILLE'
$ws.Range("J66").Value = '03 21 21 23 88'
$ws.Range("A67").Value = 66
$ws.Range("B67").Value = 174025900000066
$ws.Range("C67").Value = 'TRUSSANT'
$ws.Range("D67").Value = 'Pierre'
$ws.Range("E67").Value = 'h'
$ws.Range("F67").Value = 27085
$ws.Range("F67").PasteSpecial(-4122)
$ws.Range("G67").Value = '25 Boulevard Excellemans'
$ws.Range("H67").Value = 59000
$ws.Range("I67").Value = 'LILLE'
$ws.Range("J67").Value = '03 21 66 44 77'
$ws.Range("A68").Value = 67
$ws.Range("B68").Value = 285065900000067
$ws.Range("C68").Value = 'GIBOURI'
$ws.Range("D68").Value = 'Alicia'
$ws.Range("E68").Value = 'f'
$ws.Range("F68").Value = 31211
$ws.Range("F68").PasteSpecial(-4122)
$ws.Range("G68").Value = '30 Rue Solférino'
$ws.Range("H68").Value = 59000
$ws.Range("I68").Value = 'LILLE'
$ws.Range("J68").Value = '03 55 55 88 88'
$ws.Range("A69").Value = 68
$ws.Range("B69").Value = 266106200000068
$ws.Range("C69").Value = 'TRANCHANT'
$ws.Range("D69").Value = 'Florence'
$ws.Range("E69").Value = 'f'
$ws.Range("F69").Value = 24390
$ws.Range("F69").PasteSpecial(-4122)
$ws.Range("G69").Value = '91 Rue Nationale, apt24'
$ws.Range("H69").Value = 59000
$ws.Range("I69").Value = 'LILLE'
$ws.Range("J69").Value = '03 23 98 89 77'
$ws.Range("A70").Value = 69
$ws.Range("B70").Value = 191085900000069
$ws.Range("C70").Value = 'FILLOT'
$ws.Range("D70").Value = 'Keran'
$ws.Range("E70").Value = 'h'
$ws.Range("F70").Value = 33452
$ws.Range("F70").PasteSpecial(-4122)
$ws.Range("G70").Value = '149 Rue du Général de Gaulles'
$ws.Range("H70").Value = 59000
$ws.Range("I70").Value = 'LILLE'
$ws.Range("J70").Value = '03 45 65 56 65'
$ws.Range("A71").Value = 70
$ws.Range("B71").Value = 284095900000070
$ws.Range("C71").Value = 'DELAUME'
$ws.Range("D71").Value = 'Auriane'
$ws.Range("E71").Value = 'f'
$ws.Range("F71").Value = 30951
$ws.Range("F71").PasteSpecial(-4122)
$ws.Range("G71").Value = '2 Rue de la Monnaie'
$ws.Range("H71").Value = 59000
$ws.Range("I71").Value = 'LILLE'
$ws.Range("J71").Value = '03 21 11 44 77'

# Etat de la vue (selection / defilement) pour refleter la position finale de l edition
$ws.Range("B71").Select()
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "done"
